# Added "Max-Min HR Window 10 Sec Stats" sheet with a stress-window
# calculator (median/outlier based bucketing of a capped bpm delta).

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the tab strip ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Max-Min HR Window 10 Sec Stats"

# --- headers (write order matches the authored shared-string table) -------
$ws.Range("A1").Value = "Max-Min HR Window Calculation"
$ws.Range("C1").Value = "Median"
$ws.Range("D1").Value = "No Stress"
$ws.Range("E1").Value = "Medium Stress"
$ws.Range("F1").Value = "High Stress"
$ws.Range("B1").Value = "Max-Min HR Window Calculation Filtered Outliers"

# --- raw Max-Min HR window samples (column A) ------------------------------
$rawValues = @(4,4,5,4,13,12,6,7,2,3,9,13,12,6,7,3,7,2,3,17,20,25,8,3,4)
for ($i = 0; $i -lt $rawValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $rawValues[$i]
}

# --- filtered-outlier column (B): cap raw value at 10 ----------------------
$ws.Range("B2").Formula = "=IF(A2<10, A2, 10)"
$ws.Range("B3:B26").Formula = "=IF(A3<10, A3, 10)"

# --- stats row (row 2 only) -------------------------------------------------
$ws.Range("C2").Formula = "=MEDIAN(B2:B26)"
$ws.Range("D2").Formula = "=COUNTIF(B2:B26, ""<=""&C2)/COUNT(B2:B26)"
$ws.Range("E2").Formula = "=COUNTIF(B2:B26, "">""&C2)/COUNT(B2:B26)-F2"
$ws.Range("F2").Formula = "=COUNTIF(B2:B26, "">=10"")/COUNT(B2:B26)"

# --- scratch column H --------------------------------------------------
$ws.Cells.Item(2, 8).Value = 2
$ws.Cells.Item(3, 8).Value = 2
$ws.Cells.Item(12, 8).Value = 4
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(15, 8).Value = 5
$ws.Cells.Item(16, 8).Value = 6
$ws.Cells.Item(17, 8).Value = 7
$ws.Cells.Item(18, 8).Value = 10
$ws.Cells.Item(19, 8).Value = 10

# --- column widths (best effort match to the authored layout) --------------
$ws.Columns.Item(1).ColumnWidth = 26.92
$ws.Columns.Item(2).ColumnWidth = 32.25
$ws.Columns.Item(3).ColumnWidth = 12.42
$ws.Columns.Item(4).ColumnWidth = 9.25
$ws.Columns.Item(5).ColumnWidth = 11.92
$ws.Columns.Item(6).ColumnWidth = 9.92

# make the new sheet the active tab (mirrors tabSelected move + activeTab=3)
$ws.Activate()

# --- selection / view state --------------------------------------------
$ws.Range("H4").Select()
